# Regenerate the s_vals data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) for rows 2-22 with the
# newly-computed values, and recomputes G (sum) as B+C+D+E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697),
    @(3, 0.1190320826869504, 0.04071648406533734, 3.537761648806719, 0.4942365360607697),
    @(4, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(5, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(6, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697),
    @(7, 1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697),
    @(8, 0.1190320826869504, 0.306821227259698, 22.3905356188092, 10.19245300693656),
    @(9, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(10, 0.1190320826869504, 0.04071648406533734, 3.537761648806719, 0.4942365360607697),
    @(11, 0.04271373187048222, 0.306821227259698, 3.537761648806719, 0.4942365360607697),
    @(12, 3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697),
    @(13, 3.286832544864788, 1.655778082260271, 22.3905356188092, 10.19245300693656),
    @(14, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697),
    @(15, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(16, 1.455362044514542, 0.306821227259698, 0.7527432677738641, 0.4942365360607697),
    @(17, 0.1190320826869504, 0.04071648406533734, 3.537761648806719, 0.4942365360607697),
    @(18, 3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(19, 1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697),
    @(20, 0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(21, 1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697),
    @(22, 1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697)
)

foreach ($row in $data) {
    $r = $row[0]
    $b = $row[1]
    $c = $row[2]
    $d = $row[3]
    $e = $row[4]

    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 7).Value = $b + $c + $d + $e
}
